$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (which only held the "5963230 - Leandro Goncalves de Aguiar"
# text in columns B/C) is removed entirely; every row below it shifts up by
# one (row heights / spans follow automatically with a native row delete).
$ws.Rows.Item(13).Delete()

# After the shift, several rows' B/C (value + red "modified" copy) text is
# replaced with different content than what simply shifted into place.
$ws.Range("B10").Value = '5963230 - Leandro Gonçalves de Aguiar'
$ws.Range("C10").Value = '5963230 - Leandro Gonçalves de Aguiar'

$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

$ws.Range("B15").Value = '01/01/2012'
$ws.Range("C15").Value = '01/01/2012'

$ws.Range("B18").Value = '5963230 - Leandro Gonçalves de Aguiar'
$ws.Range("C18").Value = '5963230 - Leandro Gonçalves de Aguiar'

$ws.Range("B19").Value = 'Duas provas escritas e eventual apresentação de trabalho.'
$ws.Range("C19").Value = 'Duas provas escritas e eventual apresentação de trabalho.'

$ws.Range("B20").Value = 'Nota(N) = 50% Prova P1 + 50% Prova P2. Os pesos poderão ser redefinidos caso seja incorporada nota de trabalho.'
$ws.Range("C20").Value = 'Nota(N) = 50% Prova P1 + 50% Prova P2. Os pesos poderão ser redefinidos caso seja incorporada nota de trabalho.'

$ws.Range("B21").Value = 'Média Final = (N + Prova Recuperação)/2'
$ws.Range("C21").Value = 'Média Final = (N + Prova Recuperação)/2'
